$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loan Quotes")

$ws.Range("A9").Value = "TEST"
$ws.Range("B9").Value = 100000
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 252000
$ws.Range("E9").Value = 22
$ws.Range("F9").Value = "No"

$ws.Range("G9").Select()
